$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-25 Wednesday" "2025-06-26 Thursday"

Replace-Text "472×4=1888" "769×3=2307"
Replace-Text "682×8=5456" "916×4=3664"
Replace-Text "899×7=6293" "443×7=3101"
Replace-Text "892×6=5352" "975×4=3900"
Replace-Text "678×5=3390" "983×5=4915"

Replace-Text "699×9=6291" "585×5=2925"
Replace-Text "383×6=2298" "792×2=1584"
Replace-Text "939×2=1878" "912×3=2736"
Replace-Text "740×7=5180" "832×4=3328"
Replace-Text "125×8=1000" "163×2=326"

Replace-Text "434×4=1736" "926×5=4630"
Replace-Text "546×4=2184" "305×3=915"
Replace-Text "267×7=1869" "391×3=1173"
Replace-Text "864×5=4320" "650×6=3900"
Replace-Text "582×5=2910" "275×2=550"

Replace-Text "294×7=2058" "438×5=2190"
Replace-Text "300×4=1200" "187×3=561"
Replace-Text "245×5=1225" "211×3=633"
Replace-Text "378×4=1512" "935×5=4675"
Replace-Text "982×2=1964" "256×9=2304"

Replace-Text "741×9=6669" "985×5=4925"
Replace-Text "612×6=3672" "220×4=880"
Replace-Text "873×9=7857" "652×2=1304"
Replace-Text "410×9=3690" "395×9=3555"
Replace-Text "599×2=1198" "682×8=5456"
